$d = $word.ActiveDocument

function Get-ParaRangeByText($fullText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($fullText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Could not find paragraph text: " + $fullText)
    }
    $para = $rng.Paragraphs(1)
    return $para.Range
}

function Wrap-Pkg($bodyInner) {
    $body = '<w:body xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyInner + '</w:body>'
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`n" +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + "`n" +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + "`n" +
           '<pkg:xmlData>' + "`n" +
           $body + "`n" +
           '</pkg:xmlData>' + "`n" +
           '</pkg:part>' + "`n" +
           '</pkg:package>'
    return $pkg
}

# ---------------------------------------------------------------------------
# 1) "Кориснику се приказује избор неоткључаних жанрова."
#    -> "Корисник бира ставку тренинг режим из главног менија."
# ---------------------------------------------------------------------------
$r1 = Get-ParaRangeByText("Кориснику се приказује избор неоткључаних жанрова.")
$r1b = $d.Range($r1.Start, $r1.End - 1)
$r1b.Text = "Корисник бира ставку тренинг режим из главног менија."

# ---------------------------------------------------------------------------
# 2) "Селекцијом неких од жанрова прелази на следећи прозор. ... извођача" + "."
#    -> "Кориснику се приказује избор " + "свих постојећих жанрова, ... откључати."
# ---------------------------------------------------------------------------
$r2 = Get-ParaRangeByText("Селекцијом неких од жанрова прелази на следећи прозор. Корисник добија десет питања од којих шест морају да буду тачна. Питања се састоје од погађања песама или извођача.")
$r2b = $d.Range($r2.Start, $r2.End - 1)
$inner2 = '<w:p>' +
    '<w:r><w:rPr><w:lang w:val="sr-Cyrl-RS"/></w:rPr><w:t xml:space="preserve">Кориснику се приказује избор </w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="sr-Cyrl-RS"/></w:rPr><w:t>свих постојећих жанрова, само поред оних који нису откључани стоји индикација да се могу откључати.</w:t></w:r>' +
    '</w:p>'
$r2b.InsertXML((Wrap-Pkg $inner2))

# ---------------------------------------------------------------------------
# 3) "Песме се бирају на основу случајно изабраних нумера из свих плејлисти тог жанра."
#    -> split into several runs describing genre selection (locked genres)
# ---------------------------------------------------------------------------
$r3 = Get-ParaRangeByText("Песме се бирају на основу случајно изабраних нумера из свих плејлисти тог жанра.")
$r3b = $d.Range($r3.Start, $r3.End - 1)
$inner3 = '<w:p>' +
    '<w:r><w:rPr><w:lang w:val="sr-Cyrl-RS"/></w:rPr><w:t>Селекцијом нек</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="sr-Cyrl-RS"/></w:rPr><w:t>ог</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="sr-Cyrl-RS"/></w:rPr><w:t xml:space="preserve"> од</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="sr-Cyrl-RS"/></w:rPr><w:t>закључаних</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="sr-Cyrl-RS"/></w:rPr><w:t xml:space="preserve"> жанрова прелази на следећи прозор. Корисник добија десет питања од којих шест морају да буду тачна. Питања се састоје од погађања песама или извођача</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="sr-Cyrl-RS"/></w:rPr><w:t>.</w:t></w:r>' +
    '</w:p>'
$r3b.InsertXML((Wrap-Pkg $inner3))

# ---------------------------------------------------------------------------
# 4) "Кориснику се појављује прозор са честитком"
#    -> "Песме се бирају на основу случајно изабраних нумера из свих плејлисти тог жанра."
# ---------------------------------------------------------------------------
$r4 = Get-ParaRangeByText("Кориснику се појављује прозор са честитком")
$r4b = $d.Range($r4.Start, $r4.End - 1)
$r4b.Text = "Песме се бирају на основу случајно изабраних нумера из свих плејлисти тог жанра."

# ---------------------------------------------------------------------------
# 5) "Корисник откључава нови жанр и добија две основне плејлисте"
#    -> "Кориснику се појављује прозор са честитком"
# ---------------------------------------------------------------------------
$r5 = Get-ParaRangeByText("Корисник откључава нови жанр и добија две основне плејлисте")
$r5b = $d.Range($r5.Start, $r5.End - 1)
$r5b.Text = "Кориснику се појављује прозор са честитком"

# ---------------------------------------------------------------------------
# 6) The paragraph with "откључава могућност такмичења у том жанру и добија " + ...
#    becomes TWO paragraphs:
#      a) "Корисник откључава нови жанр и добија две основне плејлисте"
#      b) (new paragraph) lastRenderedPageBreak + the original run sequence
# ---------------------------------------------------------------------------
$r6 = Get-ParaRangeByText("откључава могућност такмичења у том жанру и добија одређени број поена на коначној ранг листи као награду")
$r6b = $d.Range($r6.Start, $r6.End - 1)
$r6b.Text = "Корисник откључава нови жанр и добија две основне плејлисте"

# insert a brand-new list paragraph right after this one, carrying the
# previous content (now with a lastRenderedPageBreak marker on the first run)
$r6After = $d.Range($r6b.End, $r6b.End)
$inner6b = '<w:p>' +
    '<w:pPr>' +
        '<w:pStyle w:val="ListParagraph"/>' +
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="16"/></w:numPr>' +
        '<w:ind w:left="993"/>' +
        '<w:rPr><w:lang w:val="sr-Cyrl-RS"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="sr-Cyrl-RS"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">откључава могућност такмичења у том жанру и добија </w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="sr-Cyrl-RS"/></w:rPr><w:t xml:space="preserve">одређени број поена </w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="sr-Cyrl-RS"/></w:rPr><w:t xml:space="preserve">на коначној </w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="sr-Cyrl-RS"/></w:rPr><w:t>ранг листи</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="sr-Cyrl-RS"/></w:rPr><w:t xml:space="preserve"> као награду</w:t></w:r>' +
    '</w:p>'
$r6After.InsertXML((Wrap-Pkg $inner6b))

# ---------------------------------------------------------------------------
# 7) Remove the (now duplicate) lastRenderedPageBreak from the heading
#    "Корисник не испуњава услов за откључавање новог жанра"
# ---------------------------------------------------------------------------
$r7 = Get-ParaRangeByText("Корисник не испуњава услов за откључавање новог жанра")
$inner7 = '<w:p>' +
    '<w:pPr>' +
        '<w:pStyle w:val="Heading3"/>' +
        '<w:rPr><w:b w:val="0"/><w:i/><w:lang w:val="sr-Cyrl-RS"/></w:rPr>' +
    '</w:pPr>' +
    '<w:bookmarkStart w:id="9" w:name="_Toc67766891"/>' +
    '<w:r><w:rPr><w:b w:val="0"/><w:i/><w:lang w:val="sr-Cyrl-RS"/></w:rPr><w:t>Корисник не испуњава услов за откључавање новог жанра</w:t></w:r>' +
    '<w:bookmarkEnd w:id="9"/>' +
    '</w:p>'
$r7.InsertXML((Wrap-Pkg $inner7))
